$wb = $excel.ActiveWorkbook

# Update "想去人数" (want-to-go count) values on both the "展览" and
# "全部类型" worksheets, which contain the same duplicated data rows.
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 92
    $ws.Range("F3").Value = 304
}
